# Auto-generated edit script for violent-crime-full-year.xlsx
# Commit: Add data for 2023-12-05
# Updates column J (year 2023) totals across 47 worksheets (158 cells total)
# to reflect newly added data for 2023-12-05.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Citywide Totals")
$ws.Cells.Item(2, 10).Value = 7121   # J2: 7104 -> 7121
$ws.Cells.Item(3, 10).Value = 7521   # J3: 7506 -> 7521
$ws.Cells.Item(4, 10).Value = 1642   # J4: 1635 -> 1642
$ws.Cells.Item(5, 10).Value = 587   # J5: 585 -> 587
$ws.Cells.Item(6, 10).Value = 10225   # J6: 10183 -> 10225
$ws.Cells.Item(7, 10).Value = 27096   # J7: 27013 -> 27096

$ws = $wb.Worksheets.Item("By Neighborhood")
$ws.Cells.Item(6, 10).Value = 209   # J6: 208 -> 209
$ws.Cells.Item(7, 10).Value = 774   # J7: 772 -> 774
$ws.Cells.Item(8, 10).Value = 1701   # J8: 1690 -> 1701
$ws.Cells.Item(9, 10).Value = 139   # J9: 138 -> 139
$ws.Cells.Item(11, 10).Value = 485   # J11: 482 -> 485
$ws.Cells.Item(16, 10).Value = 108   # J16: 105 -> 108
$ws.Cells.Item(19, 10).Value = 784   # J19: 782 -> 784
$ws.Cells.Item(20, 10).Value = 570   # J20: 568 -> 570
$ws.Cells.Item(22, 10).Value = 64   # J22: 63 -> 64
$ws.Cells.Item(24, 10).Value = 91   # J24: 88 -> 91
$ws.Cells.Item(27, 10).Value = 166   # J27: 164 -> 166
$ws.Cells.Item(29, 10).Value = 1446   # J29: 1444 -> 1446
$ws.Cells.Item(30, 10).Value = 94   # J30: 93 -> 94
$ws.Cells.Item(31, 10).Value = 279   # J31: 277 -> 279
$ws.Cells.Item(33, 10).Value = 1231   # J33: 1228 -> 1231
$ws.Cells.Item(34, 10).Value = 126   # J34: 125 -> 126
$ws.Cells.Item(36, 10).Value = 365   # J36: 364 -> 365
$ws.Cells.Item(37, 10).Value = 831   # J37: 829 -> 831
$ws.Cells.Item(39, 10).Value = 20   # J39: 19 -> 20
$ws.Cells.Item(40, 10).Value = 61   # J40: 60 -> 61
$ws.Cells.Item(42, 10).Value = 1161   # J42: 1160 -> 1161
$ws.Cells.Item(45, 10).Value = 39   # J45: 38 -> 39
$ws.Cells.Item(48, 10).Value = 307   # J48: 305 -> 307
$ws.Cells.Item(51, 10).Value = 336   # J51: 334 -> 336
$ws.Cells.Item(52, 10).Value = 691   # J52: 688 -> 691
$ws.Cells.Item(54, 10).Value = 533   # J54: 530 -> 533
$ws.Cells.Item(55, 10).Value = 429   # J55: 428 -> 429
$ws.Cells.Item(63, 10).Value = 84   # J63: 85 -> 84
$ws.Cells.Item(66, 10).Value = 82   # J66: 81 -> 82
$ws.Cells.Item(67, 10).Value = 1010   # J67: 1007 -> 1010
$ws.Cells.Item(70, 10).Value = 38   # J70: 37 -> 38
$ws.Cells.Item(76, 10).Value = 387   # J76: 385 -> 387
$ws.Cells.Item(80, 10).Value = 49   # J80: 48 -> 49
$ws.Cells.Item(83, 10).Value = 544   # J83: 543 -> 544
$ws.Cells.Item(84, 10).Value = 223   # J84: 222 -> 223
$ws.Cells.Item(85, 10).Value = 1115   # J85: 1111 -> 1115
$ws.Cells.Item(86, 10).Value = 168   # J86: 167 -> 168
$ws.Cells.Item(88, 10).Value = 288   # J88: 287 -> 288
$ws.Cells.Item(92, 10).Value = 89   # J92: 86 -> 89
$ws.Cells.Item(93, 10).Value = 115   # J93: 114 -> 115
$ws.Cells.Item(94, 10).Value = 296   # J94: 293 -> 296
$ws.Cells.Item(96, 10).Value = 298   # J96: 297 -> 298
$ws.Cells.Item(98, 10).Value = 201   # J98: 200 -> 201
$ws.Cells.Item(99, 10).Value = 416   # J99: 415 -> 416
$ws.Cells.Item(101, 10).Value = 27096   # J101: 27013 -> 27096

$ws = $wb.Worksheets.Item("West Ridge")
$ws.Cells.Item(2, 10).Value = 90   # J2: 89 -> 90
$ws.Cells.Item(7, 10).Value = 298   # J7: 297 -> 298

$ws = $wb.Worksheets.Item("Auburn Gresham")
$ws.Cells.Item(3, 10).Value = 234   # J3: 233 -> 234
$ws.Cells.Item(4, 10).Value = 32   # J4: 31 -> 32
$ws.Cells.Item(7, 10).Value = 774   # J7: 772 -> 774

$ws = $wb.Worksheets.Item("Belmont Cragin")
$ws.Cells.Item(2, 10).Value = 138   # J2: 135 -> 138
$ws.Cells.Item(7, 10).Value = 485   # J7: 482 -> 485

$ws = $wb.Worksheets.Item("South Shore")
$ws.Cells.Item(3, 10).Value = 403   # J3: 402 -> 403
$ws.Cells.Item(6, 10).Value = 317   # J6: 314 -> 317
$ws.Cells.Item(7, 10).Value = 1115   # J7: 1111 -> 1115

$ws = $wb.Worksheets.Item("Little Village")
$ws.Cells.Item(4, 10).Value = 27   # J4: 26 -> 27
$ws.Cells.Item(5, 10).Value = 12   # J5: 11 -> 12
$ws.Cells.Item(6, 10).Value = 296   # J6: 295 -> 296
$ws.Cells.Item(7, 10).Value = 691   # J7: 688 -> 691

$ws = $wb.Worksheets.Item("Austin")
$ws.Cells.Item(2, 10).Value = 449   # J2: 447 -> 449
$ws.Cells.Item(6, 10).Value = 619   # J6: 610 -> 619
$ws.Cells.Item(7, 10).Value = 1701   # J7: 1690 -> 1701

$ws = $wb.Worksheets.Item("South Chicago")
$ws.Cells.Item(2, 10).Value = 158   # J2: 157 -> 158
$ws.Cells.Item(7, 10).Value = 544   # J7: 543 -> 544

$ws = $wb.Worksheets.Item("Garfield Park")
$ws.Cells.Item(3, 10).Value = 408   # J3: 407 -> 408
$ws.Cells.Item(6, 10).Value = 440   # J6: 438 -> 440
$ws.Cells.Item(7, 10).Value = 1231   # J7: 1228 -> 1231

$ws = $wb.Worksheets.Item("Grand Crossing")
$ws.Cells.Item(6, 10).Value = 240   # J6: 238 -> 240
$ws.Cells.Item(7, 10).Value = 831   # J7: 829 -> 831

$ws = $wb.Worksheets.Item("Woodlawn")
$ws.Cells.Item(3, 10).Value = 165   # J3: 164 -> 165
$ws.Cells.Item(7, 10).Value = 416   # J7: 415 -> 416

$ws = $wb.Worksheets.Item("Fuller Park")
$ws.Cells.Item(2, 10).Value = 33   # J2: 32 -> 33
$ws.Cells.Item(7, 10).Value = 94   # J7: 93 -> 94

$ws = $wb.Worksheets.Item("Gage Park")
$ws.Cells.Item(6, 10).Value = 98   # J6: 96 -> 98
$ws.Cells.Item(7, 10).Value = 279   # J7: 277 -> 279

$ws = $wb.Worksheets.Item("North Lawndale")
$ws.Cells.Item(2, 10).Value = 257   # J2: 256 -> 257
$ws.Cells.Item(4, 10).Value = 68   # J4: 67 -> 68
$ws.Cells.Item(6, 10).Value = 281   # J6: 280 -> 281
$ws.Cells.Item(7, 10).Value = 1010   # J7: 1007 -> 1010

$ws = $wb.Worksheets.Item("South Deering")
$ws.Cells.Item(6, 10).Value = 73   # J6: 72 -> 73
$ws.Cells.Item(7, 10).Value = 223   # J7: 222 -> 223

$ws = $wb.Worksheets.Item("Loop")
$ws.Cells.Item(2, 10).Value = 133   # J2: 132 -> 133
$ws.Cells.Item(4, 10).Value = 42   # J4: 41 -> 42
$ws.Cells.Item(6, 10).Value = 247   # J6: 246 -> 247
$ws.Cells.Item(7, 10).Value = 533   # J7: 530 -> 533

$ws = $wb.Worksheets.Item("Englewood")
$ws.Cells.Item(3, 10).Value = 512   # J3: 510 -> 512
$ws.Cells.Item(7, 10).Value = 1446   # J7: 1444 -> 1446

$ws = $wb.Worksheets.Item("Lake View")
$ws.Cells.Item(2, 10).Value = 51   # J2: 50 -> 51
$ws.Cells.Item(3, 10).Value = 59   # J3: 58 -> 59
$ws.Cells.Item(7, 10).Value = 307   # J7: 305 -> 307

$ws = $wb.Worksheets.Item("Chatham")
$ws.Cells.Item(3, 10).Value = 224   # J3: 223 -> 224
$ws.Cells.Item(6, 10).Value = 302   # J6: 301 -> 302
$ws.Cells.Item(7, 10).Value = 784   # J7: 782 -> 784

$ws = $wb.Worksheets.Item("River North")
$ws.Cells.Item(2, 10).Value = 69   # J2: 68 -> 69
$ws.Cells.Item(6, 10).Value = 204   # J6: 203 -> 204
$ws.Cells.Item(7, 10).Value = 387   # J7: 385 -> 387

$ws = $wb.Worksheets.Item("Ashburn")
$ws.Cells.Item(3, 10).Value = 49   # J3: 48 -> 49
$ws.Cells.Item(7, 10).Value = 209   # J7: 208 -> 209

$ws = $wb.Worksheets.Item("Humboldt Park")
$ws.Cells.Item(4, 10).Value = 49   # J4: 48 -> 49
$ws.Cells.Item(7, 10).Value = 1161   # J7: 1160 -> 1161

$ws = $wb.Worksheets.Item("Rogers Park")
$ws.Cells.Item(3, 10).Value = 96   # J3: 97 -> 96
$ws.Cells.Item(6, 10).Value = 95   # J6: 94 -> 95

$ws = $wb.Worksheets.Item("Lower West Side")
$ws.Cells.Item(6, 10).Value = 244   # J6: 243 -> 244
$ws.Cells.Item(7, 10).Value = 429   # J7: 428 -> 429

$ws = $wb.Worksheets.Item("Dunning")
$ws.Cells.Item(3, 10).Value = 23   # J3: 22 -> 23
$ws.Cells.Item(4, 10).Value = 14   # J4: 13 -> 14
$ws.Cells.Item(6, 10).Value = 24   # J6: 23 -> 24
$ws.Cells.Item(7, 10).Value = 91   # J7: 88 -> 91

$ws = $wb.Worksheets.Item("Roseland")
$ws.Cells.Item(2, 10).Value = 211   # J2: 210 -> 211
$ws.Cells.Item(4, 10).Value = 43   # J4: 44 -> 43

$ws = $wb.Worksheets.Item("Chicago Lawn")
$ws.Cells.Item(3, 10).Value = 192   # J3: 190 -> 192
$ws.Cells.Item(7, 10).Value = 570   # J7: 568 -> 570

$ws = $wb.Worksheets.Item("Grand Boulevard")
$ws.Cells.Item(3, 10).Value = 117   # J3: 116 -> 117
$ws.Cells.Item(7, 10).Value = 365   # J7: 364 -> 365

$ws = $wb.Worksheets.Item("West Lawn")
$ws.Cells.Item(2, 10).Value = 33   # J2: 32 -> 33
$ws.Cells.Item(7, 10).Value = 115   # J7: 114 -> 115

$ws = $wb.Worksheets.Item("Garfield Ridge")
$ws.Cells.Item(3, 10).Value = 33   # J3: 32 -> 33
$ws.Cells.Item(7, 10).Value = 126   # J7: 125 -> 126

$ws = $wb.Worksheets.Item("West Loop")
$ws.Cells.Item(6, 10).Value = 156   # J6: 153 -> 156
$ws.Cells.Item(7, 10).Value = 296   # J7: 293 -> 296

$ws = $wb.Worksheets.Item("Wicker Park")
$ws.Cells.Item(6, 10).Value = 131   # J6: 130 -> 131
$ws.Cells.Item(7, 10).Value = 201   # J7: 200 -> 201

$ws = $wb.Worksheets.Item("Greektown")
$ws.Cells.Item(4, 10).Value = 3   # J4: 2 -> 3
$ws.Cells.Item(6, 10).Value = 20   # J6: 19 -> 20

$ws = $wb.Worksheets.Item("North Center")
$ws.Cells.Item(3, 10).Value = 14   # J3: 13 -> 14
$ws.Cells.Item(7, 10).Value = 82   # J7: 81 -> 82

$ws = $wb.Worksheets.Item("Avalon Park")
$ws.Cells.Item(3, 10).Value = 44   # J3: 43 -> 44
$ws.Cells.Item(7, 10).Value = 139   # J7: 138 -> 139

$ws = $wb.Worksheets.Item("West Elsdon")
$ws.Cells.Item(6, 10).Value = 30   # J6: 27 -> 30
$ws.Cells.Item(7, 10).Value = 89   # J7: 86 -> 89

$ws = $wb.Worksheets.Item("O'Hare")
$ws.Cells.Item(2, 10).Value = 17   # J2: 16 -> 17
$ws.Cells.Item(7, 10).Value = 38   # J7: 37 -> 38

$ws = $wb.Worksheets.Item("United Center")
$ws.Cells.Item(3, 10).Value = 68   # J3: 67 -> 68
$ws.Cells.Item(7, 10).Value = 288   # J7: 287 -> 288

$ws = $wb.Worksheets.Item("Edgewater")
$ws.Cells.Item(2, 10).Value = 43   # J2: 42 -> 43
$ws.Cells.Item(4, 10).Value = 21   # J4: 20 -> 21
$ws.Cells.Item(7, 10).Value = 166   # J7: 164 -> 166

$ws = $wb.Worksheets.Item("Streeterville")
$ws.Cells.Item(4, 10).Value = 90   # J4: 89 -> 90
$ws.Cells.Item(7, 10).Value = 168   # J7: 167 -> 168

$ws = $wb.Worksheets.Item("Little Italy, UIC")
$ws.Cells.Item(4, 10).Value = 29   # J4: 28 -> 29
$ws.Cells.Item(6, 10).Value = 137   # J6: 136 -> 137
$ws.Cells.Item(7, 10).Value = 336   # J7: 334 -> 336

$ws = $wb.Worksheets.Item("Clearing")
$ws.Cells.Item(6, 10).Value = 15   # J6: 14 -> 15
$ws.Cells.Item(7, 10).Value = 64   # J7: 63 -> 64

$ws = $wb.Worksheets.Item("Jackson Park")
$ws.Cells.Item(6, 10).Value = 15   # J6: 14 -> 15
$ws.Cells.Item(7, 10).Value = 39   # J7: 38 -> 39

$ws = $wb.Worksheets.Item("Rush & Division")
$ws.Cells.Item(6, 10).Value = 27   # J6: 26 -> 27
$ws.Cells.Item(7, 10).Value = 49   # J7: 48 -> 49

$ws = $wb.Worksheets.Item("Hegewisch")
$ws.Cells.Item(3, 10).Value = 21   # J3: 20 -> 21
$ws.Cells.Item(7, 10).Value = 61   # J7: 60 -> 61

$ws = $wb.Worksheets.Item("Bucktown")
$ws.Cells.Item(6, 10).Value = 84   # J6: 81 -> 84
$ws.Cells.Item(7, 10).Value = 108   # J7: 105 -> 108

Write-Host "Updated 158 cells across 47 worksheets."
